# HAWAII_2019.xlsx data-cleaning fix
#
# 1. Rename header columns (row 1) to snake_case machine-friendly names.
# 2. Title-case the Spanish connector words ("de"/"del"/"el"/"la") inside
#    state/municipality names so they read "De"/"Del"/"El"/"La".
# 3. Tiny floating point correction on D89 (recalculated percentage).
# 4. Drop the trailing footer/metadata rows (135-139) and shrink the
#    sheet's used range back down to A1:D133.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row -----------------------------------------------------
$ws.Range("A1").Value2 = "mx_state"
$ws.Range("B1").Value2 = "mx_municipality"
$ws.Range("C1").Value2 = "n_matriculas"
$ws.Range("D1").Value2 = "pct_matriculas"

# --- 2. Title-case connector words in place names -----------------------
$ws.Range("B17").Value2  = "Villa De Álvarez"
$ws.Range("A19").Value2  = "Ciudad De México"
$ws.Range("A31").Value2  = "Estado De México"
$ws.Range("B32").Value2  = "Ixtapan De La Sal"
$ws.Range("B37").Value2  = "Tlalnepantla De Baz"
$ws.Range("B48").Value2  = "Acapulco De Juárez"
$ws.Range("B49").Value2  = "Atenango Del Río"
$ws.Range("B50").Value2  = "Chilapa De Álvarez"
$ws.Range("B54").Value2  = "Tepehuacán De Guerrero"
$ws.Range("B61").Value2  = "La Manzanilla De La Paz"
$ws.Range("B62").Value2  = "Lagos De Moreno"
$ws.Range("B66").Value2  = "Tepatitlán De Morelos"
$ws.Range("B67").Value2  = "Tizapán El Alto"
$ws.Range("B69").Value2  = "Yahualica De González Gallo"
$ws.Range("B71").Value2  = "Zapotlán El Grande"
$ws.Range("B90").Value2  = "Acatlán De Pérez Figueroa"
$ws.Range("B91").Value2  = "Chalcatongo De Hidalgo"
$ws.Range("B92").Value2  = "Constancia Del Rosario"
$ws.Range("B93").Value2  = "Fresnillo De Trujano"
$ws.Range("B94").Value2  = "Ocotlán De Morelos"
$ws.Range("B95").Value2  = "Putla Villa De Guerrero"
$ws.Range("B103").Value2 = "Tlacolula De Matamoros"
$ws.Range("B108").Value2 = "Los Reyes De Juárez"
$ws.Range("B111").Value2 = "Tetela De Ocampo"
$ws.Range("B112").Value2 = "Xayacatlán De Bravo"
$ws.Range("B114").Value2 = "Amealco De Bonfil"
$ws.Range("B123").Value2 = "Ignacio De La Llave"

# --- 3. Recalculated percentage value -----------------------------------
$ws.Range("D89").Value2 = 0.09659090909090907

# --- 4. Remove trailing footer/metadata rows ----------------------------
$ws.Rows("135:139").Delete()
